$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (e.g. "31.007.91", "1.001", "0.000007816") are kept as literal text,
# matching the workbook's original inlineStr cell type.
$targetRefs = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "B18", "C18", "D18", "E18", "B19", "C19", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($r in $targetRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range("D2").Value = "31.007.91"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.958.07"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "244.63"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.4873"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").Value = "0.2951"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "0.07014"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").Value = "19.82"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "107.68"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "1.966.73"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "0.07820"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "5.484"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "0.7019"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "281.24"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").Value = "31.025.17"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000007816"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "13.32"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "2.211.42"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "5.564"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "6.513"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "9.850"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "169.13"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").Value = "19.96"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "2.192"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "1.386"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").Value = "4.632"
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("D32").Value = "1.575"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "4.467"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "0.04925"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").Value = "0.7537"
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "0.02011"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "2.690"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "6.541"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "77.76"
$ws.Range("E42").Value = "  +8.86%  "
$ws.Range("D43").Value = "0.8997"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "0.4456"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "109.15"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "8.099"
$ws.Range("E46").Value = "  +7.36%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "1.015.23"
$ws.Range("E48").Value = "  +8.80%  "
$ws.Range("D49").Value = "9.414"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "0.1253"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "35.92"
$ws.Range("E51").Value = "  -0.94%  "
